$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("organizations")

$ws.Range("B4").Value = "Phòng Dịch vụ"
$ws.Range("C4").Value = "DV"
$ws.Range("D4").Value = "Mô tả phòng dịch vụ"

$ws.Range("B5").Value = "Phòng Tự động"
$ws.Range("C5").Value = "TD"
$ws.Range("D5").Value = "Mô tả phòng tự động"

$ws.Range("B6").Value = "Phòng Giải lao"
$ws.Range("C6").Value = "GL"
$ws.Range("D6").Value = "Mô tả giải lao"

$ws.Range("B7").Value = "Phòng Hành chính"
$ws.Range("C7").Value = "HC"
$ws.Range("D7").Value = "Mô tả phòng hành chính"

$ws.Range("E4:E7").ClearContents()

$ws.Range("C5").Select()
